# This script rewrites the conversation log (rows 2-8) into two new
# conversation blocks (rows 2-7 and rows 8-13) per main_v0.py / main_v1.py.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Content) first, top-to-bottom, rows 2-6, then 9-12, so that
#     new shared strings are minted in the same order Excel produced them in ---
$ws.Range("B2").Value = 'What is your name?'
$ws.Range("B3").Value = 'I am Cuong, your AI Assistant. How can I help you today?'
$ws.Range("B4").Value = 'Hello Cuong! I''m Minh, your AI Assistant. How can I assist you today?'
$ws.Range("B5").Value = 'Hello Minh! It seems like we have a fun mix-up here. I''m Cuong, your AI Assistant. How can I assist you today?'
$ws.Range("B6").Value = 'Hello Cuong! It looks like we’re both AI Assistants here. How can I assist you today?'

$ws.Range("B7").Value = '-------------------'
$ws.Range("B8").Value = 'What is your name?'

$ws.Range("B9").Value = 'My name is Cuong, and I''m here to assist you! How can I help you today?'
$ws.Range("B10").Value = 'Hi Cuong! I''m here to help you with any questions or tasks you have. What can I assist you with today?'
$ws.Range("B11").Value = 'Hi there! It seems like you might have mixed up our roles a bit. I''m here to assist you! If you have any questions or tasks you need help with, feel free to ask. What can I do for you today?'
$ws.Range("B12").Value = 'Thank you for the clarification, Cuong! I appreciate your willingness to assist. If you have any questions or topics you''d like to discuss, feel free to let me know. How can I help you today?'

$ws.Range("B13").Value = '-------------------'

# --- Column D (RoleA_Prompt): only the conversational rows carry a prompt,
#     the separator rows (7 and 13) have none ---
$ws.Range("D2:D6").Value = 'You are Cuong. You are AI Assistant'
$ws.Range("D8:D12").Value = 'You are Cuong. You are AI Assistant'
$ws.Range("D7").ClearContents()
$ws.Range("D13").ClearContents()

# --- Column E (RoleB_Prompt): same pattern as column D ---
$ws.Range("E2:E6").Value = 'You are Minh. You are AI Assistant'
$ws.Range("E8:E12").Value = 'You are Minh. You are AI Assistant'
$ws.Range("E7").ClearContents()
$ws.Range("E13").ClearContents()

# --- Column A (Role) ---
$ws.Range("A2").Value = 'roleA'
$ws.Range("A3").Value = 'roleA'
$ws.Range("A4").Value = 'roleB'
$ws.Range("A5").Value = 'roleA'
$ws.Range("A6").Value = 'roleB'
$ws.Range("A7").Value = 'Separator'
$ws.Range("A8").Value = 'roleB'
$ws.Range("A9").Value = 'roleA'
$ws.Range("A10").Value = 'roleB'
$ws.Range("A11").Value = 'roleA'
$ws.Range("A12").Value = 'roleB'
$ws.Range("A13").Value = 'Separator'

# --- Column C (Response_Time) ---
$ws.Range("C2").Value = 1.124949216842651
$ws.Range("C3").Value = 0.7513773441314697
$ws.Range("C4").Value = 1.053507804870605
$ws.Range("C5").Value = 0.9112639427185059
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0.7712469100952148
$ws.Range("C9").Value = 0.8190987110137939
$ws.Range("C10").Value = 1.235138177871704
$ws.Range("C11").Value = 1.640358686447144
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
